# Append a new row (row 24) of bitcoin buy data for the run on 2025-06-18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matching the existing rows that
# store their dates as plain strings rather than Excel date serials, e.g.
# A10, A12:A23). Assigning the string directly via .Value would make Excel
# auto-detect "06/18/2025" as a date and convert it to a serial number, so
# instead we build it as a text formula result in a scratch cell and paste
# only the resulting value back into A24 - this keeps the cell as genuine
# text without picking up any extra number-format/style baggage.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="06/18/2025"'
$scratch.Copy()
$ws.Range("A24").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B24").Value = 0.0004780300000000008
$ws.Range("C24").Value = 104595.9458611382
$ws.Range("D24").Value = 50
